$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Sample Size" / "Seed Value" columns from the
#     run_parameters table, shifting "Output Path" / "Version" left ---
$lo = $ws.ListObjects.Item("run_parameters")

# Move D2:E3 (Output Path / Version header + values) into B2:C3,
# overwriting the Sample Size / Seed Value columns, then clear the
# now-vacated D:E cells that used to hold those values.
$ws.Range("D2:E3").Copy($ws.Range("B2:C3")) | Out-Null
$ws.Range("D2:E3").Clear() | Out-Null

# Shrink the table down to just the two remaining columns.
$lo.Resize($ws.Range("B2:C3")) | Out-Null

# Re-assert the header text so the table's column metadata re-syncs
# with the (already correct) header cell values.
$ws.Range("B2").Value2 = "Output Path"
$ws.Range("C2").Value2 = "Version"

# --- Remove the review comment that lived on the old C2 ("Seed Value") cell ---
$ws.Range("C2").Comment.Delete() | Out-Null

# --- Update the saved cursor/selection location ---
$ws.Range("Q13").Select() | Out-Null
